# Add the new "Technology" and "Real Estate" companies that were appended
# to the Companies sheet (commit: "Technology & Real Estate").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Company Name, Industry, Beta Value, Grade, TSR
$newRows = @(
    @("TSLA",  "Technology",  2.12,                 "A", 47.620000000000005),
    @("MSFT",  "Technology",  0.92,                 "B", 16.830000000000002),
    @("META",  "Technology",  1.26,                 "A", 76.72),
    @("GOOGL", "Technology",  1.08,                 "C", 35.76),
    @("NVDA",  "Technology",  1.76,                 "C", 200.93),
    @("AAPL",  "Technology",  1.29,                 "B", 22.42),
    @("UBER",  "Technology",  1.31,                 "C", 36.22),
    @("SNOW",  "Technology",  1.37,                 "A", -17.62),
    @("APLD",  "Technology",  1.55,                 "B", 81.98),
    @("PLD",   "Real Estate", 1.01,                 "B", 17.18),
    @("AMT",   "Real Estate", 0.57999999999999996,  "C", 9.7799999999999994),
    @("EQR",   "Real Estate", 0.85,                 "C", 43.37),
    @("PSA",   "Real Estate", 0.68,                 "B", 41.61),
    @("HST",   "Real Estate", 1.1499999999999999,   "B", 16),
    @("AVB",   "Real Estate", 0.97,                 "B", 45.92),
    @("ARE",   "Real Estate", 1.04,                 "B", 20.59),
    @("INVH",  "Real Estate", 0.92,                 "A", 11.15),
    @("PEAK",  "Real Estate", 0.75,                 "B", -10.84),
    @("DLR",   "Real Estate", 0.71,                 "A", 44.269999999999996)
)

$startRow = 42
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $entry = $newRows[$i]
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("D$row").Value = $entry[3]
    $ws.Range("E$row").Value = $entry[4]
}

# Leave the selection where the author left it after entering the data.
[void]$ws.Range("F38").Select()
